# September 2020 WHO Measles caseload update for DRC
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 41 was previously blank (only had number formatting applied);
# fill it in with the new reporting-period row for September 2020,
# matching the pattern used by the preceding monthly rows.
$ws.Range("A41").Value = 43831
$ws.Range("B41").Value = 44080
$ws.Range("C41").Value = 70899
$ws.Range("D41").Value = 1317
$ws.Range("E41").Value = 1026
$ws.Range("F41").Value = "September"
$ws.Range("G41").Formula = "=C41-C40"
$ws.Range("H41").Formula = "=E41-E40"

# Append the new September summary row to the Month / Total cases / Deaths
# table at the bottom of the sheet.
$ws.Range("A67").Value = "September"
$ws.Range("B67").Value = 1650
$ws.Range("C67").Value = 39
